$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "CF11 8AZ"
$ws.Range("A5").Value = "Wales"

$ws.Range("A5").Select()
